# Extend the CITES permits table with a new "2023" column (N), matching
# the diff applied to xl/worksheets/sheet1.xml:
#   - dimension A1:M5 -> A1:N5
#   - every row's <row spans="1:13"> -> spans="1:14"
#   - row 3 gets an explicit customHeight (12.75 -> 13.5, customHeight="1")
#     and a new empty, bordered N3 cell (same style as M3)
#   - row 4 gets a new N4 cell containing 2023 (same style as M4)
#   - row 5 gets a new N5 cell containing 553 (same style as M5)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column N values (year header + data point) ---
$ws.Cells.Item(4, 14).Value = 2023
$ws.Cells.Item(5, 14).Value = 553

# Row 3's N cell stays empty but still needs the same bordered style as
# the rest of that row, so pull formatting (not values) from column M
# across rows 3-5 into the new column N.
$ws.Range("M3:M5").Copy() | Out-Null
$ws.Range("N3:N5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Row 3 switches from an implicit default height to an explicit custom
# height of 13.5 points.
$ws.Rows.Item(3).RowHeight = 13.5

# Tidy up the selection (previously a stray O4 left over from earlier
# editing) back to the top-left of the sheet.
$ws.Range("A1").Select() | Out-Null
